$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter B2:B4 values as text (so the counters are stored as text, not
# numbers) - this is what trips Excel's "number stored as text" checker.
$ws.Range("B2").Value = "1"
$ws.Range("B3").Value = "2"
$ws.Range("B4").Value = "3"

# Dismiss the resulting "number stored as text" warning for B2:B4, the way
# it's done interactively (right click -> "Ignore Error").
$errs = $ws.Range("B2:B4").Errors
$errs.Item([Microsoft.Office.Interop.Excel.XlErrorChecks]::xlNumberAsText).Ignore = $true

# Fix the accented header: "TIPO" column now holds the literal x'x' marker.
$ws.Range("C1").Value = "x'x'"

# Update the active selection to match the recorded view state
$ws.Range("G15:H16").Select() | Out-Null
